$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data of row 2 and row 4 for the columns that differ between them
# (A, B, E, F, G, H, M, Q, R). Everything else in these two rows is identical,
# so only those columns need to be exchanged.

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $cellTop = $ws.Range("$col`2")
    $cellBottom = $ws.Range("$col`4")
    $tmp = $cellTop.Value2
    $cellTop.Value2 = $cellBottom.Value2
    $cellBottom.Value2 = $tmp
}

# Column M only has a value in row 2 ("färska gnagspår"); after the edit it
# belongs to row 4, and row 2's M cell becomes empty.
$ws.Range("M4").Value2 = $ws.Range("M2").Value2
$ws.Range("M2").Value2 = $null
